# Fix Training Data Issue (#48)
# The "Date" column (BF) held a malformed literal string like "5-21-2013-14"
# (sheet name glued to the actual date) on every data row. Correct it to the
# real ISO-style date string "2014-05-21" for rows 2-31.
#
# We deliberately do NOT just assign a plain string via .Value / .Value2
# because Excel's input-parsing auto-detects "2014-05-21" as a real date,
# converting the cell to a date serial number and silently attaching a new
# date number-format style to it. Instead we write a text formula that
# evaluates to the desired literal, then convert that formula to its
# resulting value in place (copy / paste-special values-only). That mirrors
# how the original data fix was produced (by script, not manual typing) and
# keeps the cell's original style/format untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateRange = $ws.Range("BF2:BF31")
$dateRange.Formula = '="2014-05-21"'
$dateRange.Copy()
$dateRange.PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = $false
